# Applies the "Updated cryptos list" data refresh to Sheet1.
#
# The touched cells hold literal text (prices/links/percent strings are
# stored as inline strings in the source workbook, not numbers). Some of
# the new values (e.g. "0.9998", "249.32") would otherwise be auto-coerced
# to a number by Excel's normal text-to-value parsing when assigned
# through .Value. To keep them as text - matching the original cells -
# without leaving a stray NumberFormat/style behind, each cell is
# temporarily marked as Text ("@"), written, then restored to the
# "Normal" style so no formatting diff is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = '30.303.48'
$c.Style = "Normal"
$c = $ws.Cells.Item(2, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.08%  '
$c.Style = "Normal"

# Row 3
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = '1.928.01'
$c.Style = "Normal"
$c = $ws.Cells.Item(3, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.63%  '
$c.Style = "Normal"

# Row 4
$c = $ws.Cells.Item(4, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.00%  '
$c.Style = "Normal"

# Row 5
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '249.32'
$c.Style = "Normal"
$c = $ws.Cells.Item(5, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.50%  '
$c.Style = "Normal"

# Row 6
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '0.7184'
$c.Style = "Normal"
$c = $ws.Cells.Item(6, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.18%  '
$c.Style = "Normal"

# Row 7
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = '0.9998'
$c.Style = "Normal"
$c = $ws.Cells.Item(7, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.02%  '
$c.Style = "Normal"

# Row 8
$c = $ws.Cells.Item(8, 5)
$c.NumberFormat = "@"
$c.Value = '  -4.82%  '
$c.Style = "Normal"

# Row 9
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = '27.81'
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 5)
$c.NumberFormat = "@"
$c.Value = '  -3.87%  '
$c.Style = "Normal"

# Row 10
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = '0.07090'
$c.Style = "Normal"
$c = $ws.Cells.Item(10, 5)
$c.NumberFormat = "@"
$c.Value = '  -3.42%  '
$c.Style = "Normal"

# Row 11
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = '0.7888'
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 5)
$c.NumberFormat = "@"
$c.Value = '  -3.53%  '
$c.Style = "Normal"

# Row 12
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = '0.07987'
$c.Style = "Normal"
$c = $ws.Cells.Item(12, 5)
$c.NumberFormat = "@"
$c.Value = '  -2.07%  '
$c.Style = "Normal"

# Row 13
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = '1.929.81'
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.58%  '
$c.Style = "Normal"

# Row 14
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = '5.379'
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 5)
$c.NumberFormat = "@"
$c.Value = '  -2.78%  '
$c.Style = "Normal"

# Row 15
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = '94.83'
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.51%  '
$c.Style = "Normal"

# Row 16
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = '14.66'
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.66%  '
$c.Style = "Normal"

# Row 17
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = '30.296.08'
$c.Style = "Normal"
$c = $ws.Cells.Item(17, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.18%  '
$c.Style = "Normal"

# Row 18
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = '257.66'
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.12%  '
$c.Style = "Normal"

# Row 19
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = '0.000008103'
$c.Style = "Normal"
$c = $ws.Cells.Item(19, 5)
$c.NumberFormat = "@"
$c.Value = '  -2.15%  '
$c.Style = "Normal"

# Row 20
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = '5.760'
$c.Style = "Normal"
$c = $ws.Cells.Item(20, 5)
$c.NumberFormat = "@"
$c.Value = '  -2.23%  '
$c.Style = "Normal"

# Row 21
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = '2.183.70'
$c.Style = "Normal"
$c = $ws.Cells.Item(21, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.41%  '
$c.Style = "Normal"

# Row 22
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '0.9998'
$c.Style = "Normal"

# Row 23
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.Style = "Normal"
$c = $ws.Cells.Item(23, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.03%  '
$c.Style = "Normal"

# Row 24
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = '6.825'
$c.Style = "Normal"
$c = $ws.Cells.Item(24, 5)
$c.NumberFormat = "@"
$c.Value = '  -2.12%  '
$c.Style = "Normal"

# Row 25
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = '9.541'
$c.Style = "Normal"
$c = $ws.Cells.Item(25, 5)
$c.NumberFormat = "@"
$c.Value = '  -3.23%  '
$c.Style = "Normal"

# Row 26
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = '164.63'
$c.Style = "Normal"
$c = $ws.Cells.Item(26, 5)
$c.NumberFormat = "@"
$c.Value = '  +2.68%  '
$c.Style = "Normal"

# Row 27
$c = $ws.Cells.Item(27, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.85%  '
$c.Style = "Normal"

# Row 28
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = '2.272'
$c.Style = "Normal"
$c = $ws.Cells.Item(28, 5)
$c.NumberFormat = "@"
$c.Value = '  -6.84%  '
$c.Style = "Normal"

# Row 29
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = '0.1272'
$c.Style = "Normal"
$c = $ws.Cells.Item(29, 5)
$c.NumberFormat = "@"
$c.Value = '  -3.91%  '
$c.Style = "Normal"

# Row 31
$c = $ws.Cells.Item(31, 5)
$c.NumberFormat = "@"
$c.Value = '  -2.12%  '
$c.Style = "Normal"

# Row 32
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = '4.402'
$c.Style = "Normal"
$c = $ws.Cells.Item(32, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.91%  '
$c.Style = "Normal"

# Row 33
$c = $ws.Cells.Item(33, 5)
$c.NumberFormat = "@"
$c.Value = '  -2.52%  '
$c.Style = "Normal"

# Row 34
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = '0.05143'
$c.Style = "Normal"
$c = $ws.Cells.Item(34, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.88%  '
$c.Style = "Normal"

# Row 35
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = '1.265'
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.46%  '
$c.Style = "Normal"

# Row 36
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = '0.7447'
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.41%  '
$c.Style = "Normal"

# Row 37
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = '2.768'
$c.Style = "Normal"
$c = $ws.Cells.Item(37, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.10%  '
$c.Style = "Normal"

# Row 38
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = '0.01989'
$c.Style = "Normal"
$c = $ws.Cells.Item(38, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.16%  '
$c.Style = "Normal"

# Row 39
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = '2.798'
$c.Style = "Normal"
$c = $ws.Cells.Item(39, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.42%  '
$c.Style = "Normal"

# Row 40
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = '78.14'
$c.Style = "Normal"
$c = $ws.Cells.Item(40, 5)
$c.NumberFormat = "@"
$c.Value = '  -3.55%  '
$c.Style = "Normal"

# Row 41
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '6.370'
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 5)
$c.NumberFormat = "@"
$c.Value = '  -3.67%  '
$c.Style = "Normal"

# Row 42
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = '0.4513'
$c.Style = "Normal"
$c = $ws.Cells.Item(42, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.33%  '
$c.Style = "Normal"

# Row 43
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = '1.998'
$c.Style = "Normal"
$c = $ws.Cells.Item(43, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.58%  '
$c.Style = "Normal"

# Row 44
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = '0.8464'
$c.Style = "Normal"
$c = $ws.Cells.Item(44, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.29%  '
$c.Style = "Normal"

# Row 45
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = '0.9995'
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.04%  '
$c.Style = "Normal"

# Row 46
$c = $ws.Cells.Item(46, 2)
$c.NumberFormat = "@"
$c.Value = 'Quant'
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = '100.76'
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 5)
$c.NumberFormat = "@"
$c.Value = '  -2.95%  '
$c.Style = "Normal"

# Row 47
$c = $ws.Cells.Item(47, 2)
$c.NumberFormat = "@"
$c.Value = 'EnergySwap'
$c.Style = "Normal"
$c = $ws.Cells.Item(47, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c.Style = "Normal"
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = '9.805'
$c.Style = "Normal"
$c = $ws.Cells.Item(47, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.67%  '
$c.Style = "Normal"

# Row 48
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '7.435'
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.08%  '
$c.Style = "Normal"

# Row 49
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = '36.84'
$c.Style = "Normal"
$c = $ws.Cells.Item(49, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.53%  '
$c.Style = "Normal"

# Row 50
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = '954.97'
$c.Style = "Normal"
$c = $ws.Cells.Item(50, 5)
$c.NumberFormat = "@"
$c.Value = '  +8.50%  '
$c.Style = "Normal"

# Row 51
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = '0.4222'
$c.Style = "Normal"
$c = $ws.Cells.Item(51, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.52%  '
$c.Style = "Normal"
